$wb = $excel.ActiveWorkbook

# Sheet 1 (ALC)
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 2211.7856
$ws.Range("I28").Value = 915.9
$ws.Range("J28").Value = 5451.5
$ws.Range("K28").Value = 915.9
$ws.Range("L28").Value = 5451.5
$ws.Range("M28").Value = -430.9
$ws.Range("N28").Value = -6421.5
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 132
$ws.Range("H132").Value = 1735.5
$ws.Range("I132").Value = 1580.2222
$ws.Range("K132").Value = 4740.6666
$ws.Range("M132").Value = -2210.6666
# Row 133
$ws.Range("H133").Value = 68996.73
$ws.Range("J133").Value = 68996.73
$ws.Range("L133").Value = 68996.73
$ws.Range("N133").Value = -79116.73
# Row 134
$ws.Range("H134").Value = 98950
$ws.Range("J134").Value = 98950
$ws.Range("L134").Value = 98950
$ws.Range("N134").Value = -109090
# Row 137
$ws.Range("H137").Value = 692859.0600000001
$ws.Range("J137").Value = 969268.1
$ws.Range("L137").Value = 2907804.3
$ws.Range("N137").Value = -2912904.3
# Row 138
$ws.Range("H138").Value = 3025.6191
$ws.Range("I138").Value = 1403.1765
$ws.Range("J138").Value = 4128.88
$ws.Range("K138").Value = 4209.529500000001
$ws.Range("L138").Value = 12386.64
$ws.Range("M138").Value = 930.4704999999994
$ws.Range("N138").Value = -22666.64
# Row 139
$ws.Range("H139").Value = 70646.2
$ws.Range("J139").Value = 70646.2
$ws.Range("L139").Value = 70646.2
$ws.Range("N139").Value = -80926.2

# Sheet 2 (ARM)
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3578.0303
$ws.Range("I32").Value = 3125.1304
$ws.Range("K32").Value = 3125.1304
$ws.Range("M32").Value = -2838.1304
# Row 39
$ws.Range("H39").Value = 14250
$ws.Range("I39").Value = 9000
$ws.Range("J39").Value = 30000
$ws.Range("K39").Value = 9000
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = -8480
$ws.Range("N39").Value = -31040
# Row 61
$ws.Range("H61").Value = 74556.5
$ws.Range("I61").Value = 3253.818
$ws.Range("K61").Value = 3253.818
$ws.Range("M61").Value = -3041.818
# Row 74
$ws.Range("H74").Value = 5395
# Row 77
$ws.Range("H77").Value = 5395
# Row 136
$ws.Range("H136").Value = 74556.5
$ws.Range("I136").Value = 3253.818
$ws.Range("K136").Value = 9761.454000000002
$ws.Range("M136").Value = -7211.454000000002
# Row 139
$ws.Range("H139").Value = 128940.8
$ws.Range("J139").Value = 128940.8
$ws.Range("L139").Value = 128940.8
$ws.Range("N139").Value = -139220.8

# Sheet 3 (BSM)
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 52558.65
$ws.Range("J20").Value = 3374.5
$ws.Range("L20").Value = 3374.5
$ws.Range("N20").Value = -3868.5
# Row 132
$ws.Range("H132").Value = 27399.6
$ws.Range("J132").Value = 27399.6
$ws.Range("L132").Value = 27399.6
$ws.Range("N132").Value = -37519.6
# Row 134
$ws.Range("H134").Value = 4294.593
$ws.Range("I134").Value = 3211.4285
$ws.Range("J134").Value = 8085.6665
$ws.Range("K134").Value = 9634.2855
$ws.Range("L134").Value = 24256.9995
$ws.Range("M134").Value = -7099.2855
$ws.Range("N134").Value = -29326.9995
# Row 135
$ws.Range("H135").Value = 97617.86
$ws.Range("J135").Value = 97617.86
$ws.Range("L135").Value = 97617.86
$ws.Range("N135").Value = -107757.86
# Row 138
$ws.Range("H138").Value = 76664.44500000001
$ws.Range("J138").Value = 76664.44500000001
$ws.Range("L138").Value = 76664.44500000001
$ws.Range("N138").Value = -86944.44500000001

# Sheet 4 (CRP)
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2302.68
$ws.Range("I31").Value = 1574.1818
$ws.Range("J31").Value = 2875.0715
$ws.Range("K31").Value = 1574.1818
$ws.Range("L31").Value = 2875.0715
$ws.Range("M31").Value = -1279.1818
$ws.Range("N31").Value = -3465.0715
# Row 34
$ws.Range("H34").Value = 2302.68
$ws.Range("I34").Value = 1574.1818
$ws.Range("J34").Value = 2875.0715
$ws.Range("K34").Value = 1574.1818
$ws.Range("L34").Value = 2875.0715
$ws.Range("M34").Value = -1372.1818
$ws.Range("N34").Value = -3279.0715
# Row 95
$ws.Range("H95").Value = 14392.571
$ws.Range("J95").Value = 14392.571
$ws.Range("L95").Value = 14392.571
$ws.Range("N95").Value = -19884.571
# Row 108
$ws.Range("H108").Value = 66768.37
$ws.Range("J108").Value = 66768.37
$ws.Range("L108").Value = 66768.37
$ws.Range("N108").Value = -74448.37
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
# Row 111
$ws.Range("H111").Value = 73248.5
$ws.Range("J111").Value = 73248.5
$ws.Range("L111").Value = 73248.5
$ws.Range("N111").Value = -81428.5
# Row 130
$ws.Range("H130").Value = 150000
$ws.Range("J130").Value = 150000
$ws.Range("L130").Value = 150000
$ws.Range("N130").Value = -160040
# Row 138
$ws.Range("H138").Value = 52964.445
$ws.Range("J138").Value = 52964.445
$ws.Range("L138").Value = 52964.445
$ws.Range("N138").Value = -63244.445

# Sheet 5 (CUL)
$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 18473.5
$ws.Range("I32").Value = 46466.5
$ws.Range("K32").Value = 139399.5
$ws.Range("M32").Value = -139116.5
# Row 46
$ws.Range("H46").Value = 308
$ws.Range("I46").Value = 210
$ws.Range("J46").Value = 373.33334
$ws.Range("K46").Value = 630
$ws.Range("L46").Value = 1120.00002
$ws.Range("M46").Value = -539
$ws.Range("N46").Value = -1302.00002
# Row 68
$ws.Range("H68").Value = 1880.75
$ws.Range("J68").Value = 1930.1538
$ws.Range("L68").Value = 5790.4614
$ws.Range("N68").Value = -7412.4614
# Row 71
$ws.Range("H71").Value = 1880.75
$ws.Range("J71").Value = 1930.1538
$ws.Range("L71").Value = 17371.3842
$ws.Range("N71").Value = -25483.3842
# Row 107
$ws.Range("H107").Value = 2028.5
$ws.Range("J107").Value = 2159
$ws.Range("L107").Value = 6477
$ws.Range("N107").Value = -10317
# Row 131
$ws.Range("H131").Value = 1495
$ws.Range("J131").Value = 2190
$ws.Range("L131").Value = 6570
$ws.Range("N131").Value = -16650
# Row 138
$ws.Range("H138").Value = 6311.385
$ws.Range("I138").Value = 5000
$ws.Range("K138").Value = 15000
$ws.Range("M138").Value = -9860

# Sheet 6 (GSM)
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3534.375
$ws.Range("I80").Value = 2817.25
$ws.Range("J80").Value = 4251.5
$ws.Range("K80").Value = 2817.25
$ws.Range("L80").Value = 4251.5
$ws.Range("M80").Value = -1819.25
$ws.Range("N80").Value = -6247.5
# Row 83
$ws.Range("H83").Value = 3534.375
$ws.Range("I83").Value = 2817.25
$ws.Range("J83").Value = 4251.5
$ws.Range("K83").Value = 14086.25
$ws.Range("L83").Value = 21257.5
$ws.Range("M83").Value = -9094.25
$ws.Range("N83").Value = -31241.5
# Row 97
$ws.Range("H97").Value = 4250
$ws.Range("I97").Value = 4333.3335
$ws.Range("K97").Value = 4333.3335
$ws.Range("M97").Value = -3837.3335
# Row 113
$ws.Range("H113").Value = 2273.9443
$ws.Range("I113").Value = 2168.6667
$ws.Range("K113").Value = 2168.6667
$ws.Range("M113").Value = 1.333299999999781
# Row 126
$ws.Range("H126").Value = 6000
$ws.Range("I126").Value = 3100
$ws.Range("K126").Value = 9300
$ws.Range("M126").Value = -6830
# Row 135
$ws.Range("H135").Value = 94996.664
$ws.Range("J135").Value = 94996.664
$ws.Range("L135").Value = 94996.664
$ws.Range("N135").Value = -105136.664
# Row 136
$ws.Range("H136").Value = 86258
$ws.Range("J136").Value = 86258
$ws.Range("L136").Value = 258774
$ws.Range("N136").Value = -263874
# Row 140
$ws.Range("H140").Value = 92424.234
$ws.Range("J140").Value = 92424.234
$ws.Range("L140").Value = 92424.234
$ws.Range("N140").Value = -102784.234

# Sheet 7 (LTW)
$ws = $wb.Worksheets.Item("LTW")
# Row 92
$ws.Range("H92").Value = 100000
$ws.Range("J92").Value = 100000
$ws.Range("L92").Value = 100000
$ws.Range("N92").Value = -104992
# Row 93
$ws.Range("H93").Value = 2488.25
$ws.Range("I93").Value = 2488.25
$ws.Range("K93").Value = 2488.25
$ws.Range("M93").Value = -1240.25
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
# Row 132
$ws.Range("H132").Value = 1855.2858
$ws.Range("I132").Value = 1855.2858
$ws.Range("K132").Value = 5565.857400000001
$ws.Range("M132").Value = -3035.857400000001

# Sheet 8 (WVR)
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1561.4642
$ws.Range("I132").Value = 1205.5652
$ws.Range("K132").Value = 3616.6956
$ws.Range("M132").Value = -1086.6956
# Row 136
$ws.Range("H136").Value = 2367.2856
$ws.Range("I136").Value = 2367.2856
$ws.Range("K136").Value = 7101.8568
$ws.Range("M136").Value = -4551.8568
